# Update the "Förändrad" (Changed) date column (C) from 2023-12-10 (45270)
# to 2023-12-11 (45271) for all data rows (rows 2 through 27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45270) {
        $cell.Value2 = 45271
    }
}
